# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Apr 27 07:06:22 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.852.74"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.884.74"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'331.90"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4620"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").Value = "'0.4103"
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "'47.58"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'0.07980"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.9976"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'21.70"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "1.896.15"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "'5.907"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "'7.044"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'89.07"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06578"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.00001027"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'17.46"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "28.967.07"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "'5.422"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'2.213"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").Value = "2.120.23"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "'157.42"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'19.69"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'2.123"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'5.424"
$ws.Range("D31").Value = "'117.64"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").Value = "'0.9774"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "'0.09359"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").Value = "'1.413"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("D35").Value = "'3.601"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "'5.285"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'0.06061"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'0.02236"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'8.336"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'1.173"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'0.5794"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").Value = "'10.14"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'0.1816"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "'1.249"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "'2.287"
$ws.Range("E46").Value = "  +10.82%  "
$ws.Range("D47").Value = "'0.5480"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'11.96"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "'1.909"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "'0.07022"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'46.84"
$ws.Range("E51").Value = "  +18.20%  "
